$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new pair of rows (Primera / Segunda quality) is inserted
# right before the current first data block for this week, pushing every
# existing row down by 2 (201-268 -> 203-270).
$ws.Rows("201:202").Insert()

# Row 201 - Apio, Americana (o), Primera
$ws.Range("A201").Value = 8
$ws.Range("B201").Value = "Terminal La Palmera de La Serena"
$ws.Range("C201").Value = "Coquimbo"
$ws.Range("D201").Value = 44524
$ws.Range("E201").Value = 4
$ws.Range("F201").Value = 100112017
$ws.Range("G201").Value = "Apio"
$ws.Range("H201").Value = "Americana (o)"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 2500
$ws.Range("K201").Value = 7000
$ws.Range("L201").Value = 8000
$ws.Range("M201").Value = 7500
$ws.Range("N201").Value = "`$/docena de matas"
$ws.Range("O201").Value = "Provincia del Elquí"
$ws.Range("P201").Value = 1250
$ws.Range("Q201").Value = 6
$ws.Range("R201").Value = "Hortaliza"

# Row 202 - Apio, Americana (o), Segunda
$ws.Range("A202").Value = 8
$ws.Range("B202").Value = "Terminal La Palmera de La Serena"
$ws.Range("C202").Value = "Coquimbo"
$ws.Range("D202").Value = 44524
$ws.Range("E202").Value = 4
$ws.Range("F202").Value = 100112017
$ws.Range("G202").Value = "Apio"
$ws.Range("H202").Value = "Americana (o)"
$ws.Range("I202").Value = "Segunda"
$ws.Range("J202").Value = 1400
$ws.Range("K202").Value = 5500
$ws.Range("L202").Value = 6000
$ws.Range("M202").Value = 5750
$ws.Range("N202").Value = "`$/docena de matas"
$ws.Range("O202").Value = "Provincia del Elquí"
$ws.Range("P202").Value = 958
$ws.Range("Q202").Value = 6
$ws.Range("R202").Value = "Hortaliza"
